$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 25.567248046653923
$ws.Range("B2").Value = 1.4806419411735097
$ws.Range("C2").Value = -18.980927807839286
$ws.Range("D2").Value = 4.2866479985776529
$ws.Range("E2").Value = 408.93241489788585
$ws.Range("F2").Value = 7.7851443922479158
$ws.Range("G2").Value = 569.10835277711749
$ws.Range("H2").Value = 8.6483384226111006
$ws.Range("I2").Value = 0.058366561563422834
$ws.Range("J2").Value = 0.058366561563422834
$ws.Range("K2").Value = 10.717531205371593
$ws.Range("L2").Value = 1.531075886481656
$ws.Range("M2").Value = 0.99709377031135826
$ws.Range("N2").Value = -0.54623746562753994
